# Code divide into the different test cases
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the DealId (B2) and IdNumber (A2) values.
# Leading apostrophe keeps these as text (quote-prefixed), matching the
# existing text/number-as-text cell formatting used in the sheet.
$ws.Range("B2").Value = "'1345860"
$ws.Range("A2").Value = "'6910214183083"

# Move the active selection to A3
$ws.Range("A3").Select()
